$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 70: JNI / Intro ---
$ws.Range("A70").Value = "JNI"
$ws.Range("A71").Value = "JNA"
$ws.Range("B70").Value = "Intro"
$ws.Range("B71").Value = "Intro"
$jniBody = @"
REF: http://www.cnblogs.com/lanxuezaipiao/p/3635556.html
JNI(Java Native Interface) 有过不同语言间通信经历的一般都知道，它允许Java代码和其他语言（尤其C/C++）写的代码进行交互，只要遵守调用约定即可。首先看下JNI调用C/C++的过程，注意写程序时自下而上，调用时自上而下。Java 和 C 在開發時需互相得知對方的存在。
步骤非常的多，很麻烦，使用JNI调用.dll/.so共享库都能体会到这个痛苦的过程。如果已有一个编译好的.dll/.so文件，如果使用JNI技术调用，我们首先需要使用C语言另外写一个.dll/.so共享库，使用SUN规定的数据结构替代C语言的数据结构，调用已有的 dll/so中公布的函 数。然后再在Java中载入这个库dll/so，最后编写Java native 函数作为链接库中函数的代理。经过这些繁琐的步骤才能在Java中调用本地代码。因此，很少有Java程序员愿意编写调用dll/.so库中原生函数的java程序。这也使Java语言在客户端上乏善可陈，可以说JNI是 Java的一大弱点！
"@

$ws.Range("C70").Value = $jniBody
$jnaBody = @"
REF: http://www.cnblogs.com/lanxuezaipiao/p/3635556.html
JNA(Java Native Access) 一个开源(Github)的Java框架，是Sun公司推出的一种调用本地方法的技术，是建立在经典的JNI基础之上的一个框架。之所以说它是JNI的替代者，是因为JNA大大简化了调用本地方法的过程，使用很方便，基本上不需要脱离Java环境就可以完成。最重要的是我们不需要重写我们的动态链接库文件，而是有直接调用的API，大大简化了我们的工作量。JNA只需要我们写Java代码而不用写JNI或本地代码。
"@

$ws.Range("C71").Value = $jnaBody

# --- Row 72: JNA / Call Library ---
$ws.Range("B72").Value = "Call Library"
$ws.Range("C72").Value = "JNA can load system library like msvcrt.dll, or custom dll."

# --- Row 73: JNA / example title+ref ---
$ws.Range("B73").Value = "A complete Java - C++ connection example with mingw/g++ and JNA "
$refExample = @"
ref: http://capsis.cirad.fr/capsis/documentation/java-c_connectioncompleteexample
Super cool tutorial, I use this to finish my demo. And for code, refers to mylesieong github project: jna-demo
"@

$ws.Range("C73").Value = $refExample

# --- fill remaining column-A cells (reuse existing "JNA" string) ---
$ws.Range("A72").Value = "JNA"
$ws.Range("A73").Value = "JNA"

# --- match formatting of the rest of the table (wrap text, small font) for new rows ---
$ws.Range("A70:C73").WrapText = $true
$ws.Range("A70:A73").EntireRow.RowHeight = 33

# --- update selection / view to match the new bottom of the table ---
[void]$ws.Range("C74").Select()
$excel.ActiveWindow.ScrollRow = 64
